$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format so numeric-looking strings (e.g. "89.272.50", "0.0000294")
# are stored as text, matching the original inlineStr cell semantics.
$cellValues = [ordered]@{
    "D2" = '89.272.50'
    "E2" = '  +9.78%  '
    "D3" = '3.435.28'
    "E3" = '  +9.14%  '
    "E4" = '  +0.06%  '
    "D5" = '224.51'
    "E5" = '  +7.60%  '
    "D6" = '655.20'
    "E6" = '  +5.44%  '
    "D7" = '0.434'
    "E7" = '  +55.17%  '
    "D8" = '0.697'
    "E8" = '  +20.01%  '
    "E9" = '  -0.01%  '
    "D10" = '3.435.11'
    "E10" = '  +9.18%  '
    "D11" = '0.650'
    "E11" = '  +11.91%  '
    "D12" = '0.0000294'
    "E12" = '  +18.16%  '
    "D13" = '37.59'
    "E13" = '  +20.21%  '
    "D14" = '0.171'
    "E14" = '  +4.06%  '
    "D15" = '5.71'
    "E15" = '  +8.26%  '
    "D16" = '4.042.35'
    "E16" = '  +8.39%  '
    "D17" = '89.112.45'
    "D18" = '3.403.71'
    "E18" = '  +8.37%  '
    "D19" = '15.31'
    "E19" = '  +10.51%  '
    "B20" = 'SuiNetwork'
    "C20" = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
    "D20" = '3.17'
    "E20" = '  +0.26%  '
    "D21" = '9.84'
    "E21" = '  +10.76%  '
    "B22" = 'BitcoinCash'
    "C22" = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
    "D22" = '468.90'
    "E22" = '  +8.93%  '
    "D23" = '5.83'
    "E23" = '  +14.57%  '
    "B24" = 'NEARProtocol'
    "C24" = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
    "D24" = '5.69'
    "E24" = '  +9.21%  '
    "B25" = 'LEO'
    "C25" = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
    "D25" = '7.47'
    "E25" = '  +3.55%  '
    "D26" = '13.05'
    "E26" = '  +21.12%  '
    "B27" = 'PEPE'
    "C27" = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
    "D27" = '0.0000148'
    "E27" = '  +23.46%  '
    "B28" = 'WrappedeETH'
    "C28" = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
    "D28" = '3.585.05'
    "E28" = '  +8.18%  '
    "D29" = '82.30'
    "E29" = '  +8.03%  '
    "E30" = '  +0.06%  '
    "D31" = '0.191'
    "E31" = '  +36.17%  '
    "D32" = '9.66'
    "E32" = '  +8.44%  '
    "D33" = '598.48'
    "E33" = '  +2.46%  '
    "B34" = 'Fetch.AI'
    "C34" = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
    "D34" = '1.57'
    "E34" = '  +4.89%  '
    "B35" = 'Binance-PegBSC-USD'
    "C35" = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
    "D35" = '0.998'
    "E35" = '  -0.07%  '
    "D36" = '2.15'
    "E36" = '  +8.31%  '
    "D37" = '7.43'
    "E37" = '  +23.74%  '
    "D38" = '0.148'
    "E38" = '  -3.58%  '
    "D39" = '24.55'
    "E39" = '  +8.41%  '
    "D40" = '0.446'
    "E40" = '  +10.26%  '
    "D41" = '2.18'
    "E41" = '  +5.70%  '
    "E42" = '  +8.48%  '
    "D43" = '21.86'
    "E43" = '  +5.52%  '
    "D44" = '0.999'
    "E44" = '  +0.08%  '
    "D45" = '158.94'
    "E45" = '  +0.27%  '
    "B46" = 'Aave'
    "C46" = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
    "D46" = '193.19'
    "E46" = '  +3.77%  '
    "D47" = '1.46'
    "E47" = '  +10.73%  '
    "B48" = 'USDe'
    "C48" = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
    "D48" = '1.00'
    "E48" = '  +0.06%  '
    "D49" = '47.24'
    "E49" = '  +4.40%  '
    "D50" = '4.63'
    "E50" = '  +11.24%  '
    "D51" = '0.688'
    "E51" = '  +10.46%  '
}

foreach ($ref in $cellValues.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $cellValues[$ref]
}
